$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header shared-string in-place text edits (keep rich-text runs) ---
$ws.Range("A8").Characters(21, 2).Text = "48"
$ws.Range("C9").Characters(27, 10).Text = "11/24/2025"
$ws.Range("C9").Characters(48, 10).Text = "11/30/2025"

# --- Numeric -> text cells: force-text via leading apostrophe, then restore the
#     original (non-text) number format by pasting formats from a same-style donor ---
$ws.Range("C28").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("D28").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("E27").Value = "'***.*"
$ws.Range("E28").Value = "'***.*"
$ws.Range("F27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F27").PasteSpecial(-4122)

# --- Text -> numeric cells: paste formats from a numeric donor, then assign the number ---
$ws.Range("I14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("C18").Value = 3
$ws.Range("F31").Value = 1

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -72.727272727272
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -46.666666666666
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -41.666666666666
$ws.Range("I16").Value = 128
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -21.951219512195
$ws.Range("M16").Value = -49.407114624505
$ws.Range("N16").Value = -88.213627992633
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -42.857142857142
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 291
$ws.Range("J17").Value = 303
$ws.Range("K17").Value = -3.960396039603
$ws.Range("L17").Value = -17.796610169491
$ws.Range("M17").Value = 35.348837209302
$ws.Range("N17").Value = -13.392857142857
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 92
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = 1.098901098901
$ws.Range("L18").Value = -29.230769230769
$ws.Range("M18").Value = -72.205438066465
$ws.Range("N18").Value = -93.118922961854
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -31.25
$ws.Range("I19").Value = 315
$ws.Range("J19").Value = 331
$ws.Range("K19").Value = -4.833836858006
$ws.Range("L19").Value = -0.316455696202
$ws.Range("M19").Value = -4.255319148936
$ws.Range("N19").Value = -40.566037735849
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -11.764705882352
$ws.Range("I20").Value = 289
$ws.Range("J20").Value = 266
$ws.Range("K20").Value = 8.646616541353
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 8.239700374531
$ws.Range("N20").Value = -90.946115288220
$ws.Range("C21").Value = 20
$ws.Range("E21").Value = -23.076923076923
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = -26.086956521739
$ws.Range("I21").Value = 1134
$ws.Range("J21").Value = 1180
$ws.Range("K21").Value = -3.898305084745
$ws.Range("L21").Value = -11.128526645768
$ws.Range("M21").Value = -19.688385269121
$ws.Range("N21").Value = -82.612695492180
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 40
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -25.925925925925
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 125
$ws.Range("H24").Value = -28
$ws.Range("I24").Value = 1149
$ws.Range("J24").Value = 1320
$ws.Range("K24").Value = -12.954545454545
$ws.Range("L24").Value = -9.740769835035
$ws.Range("M24").Value = 57.181942544459
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -66.666666666666
$ws.Range("G25").Value = 86
$ws.Range("H25").Value = -55.813953488372
$ws.Range("I25").Value = 550
$ws.Range("J25").Value = 698
$ws.Range("K25").Value = -21.203438395415
$ws.Range("L25").Value = 30.641330166270
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = -6.25
$ws.Range("I26").Value = 514
$ws.Range("J26").Value = 582
$ws.Range("K26").Value = -11.683848797250
$ws.Range("L26").Value = -0.772200772200
$ws.Range("M26").Value = -11.226252158894
$ws.Range("H27").Value = -100
$ws.Range("L27").Value = -14.285714285714
$ws.Range("F28").Value = 9
$ws.Range("H28").Value = 80
$ws.Range("N29").Value = -90.625
$ws.Range("N30").Value = -88.888888888888
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 5
$ws.Range("K31").Value = -44.444444444444
$ws.Range("L31").Value = -28.571428571428
